$wb = $excel.ActiveWorkbook

# Rename sheets to unify the DataNode / DataTable naming convention.
$wb.Worksheets.Item("Property1").Name = "DataNode_1"
$wb.Worksheets.Item("Property2").Name = "DataNode_2"
$wb.Worksheets.Item("Record_Hero").Name = "DataTable_Hero"
$wb.Worksheets.Item("Record_Bag").Name = "DataTable_Bag"
$wb.Worksheets.Item("Record_CommPropertyValue").Name = "DataTable_CommPropertyValue"
$wb.Worksheets.Item("Record_Task").Name = "DataTable_Task"

# Record_Building is no longer part of the unified concept set - remove it.
$wb.Worksheets.Item("Record_Building").Delete()

# Select DataTable_Hero (formerly Record_Hero) as the active sheet/tab.
$wb.Worksheets.Item("DataTable_Hero").Activate()
